$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename shared string used by B9: "Thomas Hex" -> "Matthies Hex"
$ws.Range("B9").Value = "Matthies Hex"

# 2. Rename B4/B5 labels for newly-introduced materials "Holden" and "Rizzie Spiral"
#    (A4/A5 index values 2/3 stay as-is; row data below is fully replaced with new simulation output)
$ws.Range("B4").Value = "Holden"
$ws.Range("B5").Value = "Rizzie Spiral"

# 3. Add two new data rows (30, 31) for "Michael-CCHex" (idx 28) and "Michael-SNHex" (idx 29)
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A30").Value = 28
$ws.Range("A31").Value = 29
$ws.Range("B30").Value = "Michael-CCHex"
$ws.Range("B31").Value = "Michael-SNHex"

# 4. Rewrite the simulation data block (C4:W31) with the re-run simulation output
$arr = New-Object 'object[,]' 28,21
$arr[0,0] = 0.8297629789073429
$arr[0,1] = 1.23508634909403
$arr[0,2] = 0.9436532248701144
$arr[0,3] = 1.02082945210831
$arr[0,4] = 0.8297629789073429
$arr[0,5] = 0.960755595370459
$arr[0,6] = 0.8780909728672018
$arr[0,7] = 0.9436532248701144
$arr[0,8] = 0.9436532248701144
$arr[0,9] = 1.115468154351949
$arr[0,10] = 0.9287471941219848
$arr[0,11] = 0.9436532248701144
$arr[0,12] = 1.23508634909403
$arr[0,13] = 1.032424664000687
$arr[0,14] = 1.081916771608008
$arr[0,15] = 1.002834184290496
$arr[0,16] = 0.9978655073744527
$arr[0,17] = 1.002834184290496
$arr[0,18] = 0.9843124367483681
$arr[0,19] = 0.9761805943727173
$arr[0,20] = 0.9890492402114239
$arr[1,0] = 1.313038111058897
$arr[1,1] = 0.6309581088817445
$arr[1,2] = 1.900281771655763
$arr[1,3] = 0.6356718308312986
$arr[1,4] = 1.313038111058897
$arr[1,5] = 1.170963114224045
$arr[1,6] = 0.9910472637975943
$arr[1,7] = 1.900281771655763
$arr[1,8] = 1.900281771655763
$arr[1,9] = 0.9711712363086175
$arr[1,10] = 0.9001023562881852
$arr[1,11] = 1.900281771655763
$arr[1,12] = 0.6309581088817445
$arr[1,13] = 0.9719981099703208
$arr[1,14] = 0.7655302325849649
$arr[1,15] = 1.281425997198801
$arr[1,16] = 0.9480328587429422
$arr[1,17] = 1.281425997198801
$arr[1,18] = 1.186095086971147
$arr[1,19] = 1.32893242390807
$arr[1,20] = 1.064154224130768
$arr[2,0] = 0.8986369109390666
$arr[2,1] = 1.205658809073799
$arr[2,2] = 0.9182952508329895
$arr[2,3] = 1.029721593538399
$arr[2,4] = 0.8986369109390666
$arr[2,5] = 0.9486494071342561
$arr[2,6] = 0.9165483207142882
$arr[2,7] = 0.9182952508329895
$arr[2,8] = 0.9182952508329895
$arr[2,9] = 1.090970945582834
$arr[2,10] = 0.9384513449126813
$arr[2,11] = 0.9182952508329895
$arr[2,12] = 1.205658809073799
$arr[2,13] = 1.052147860006433
$arr[2,14] = 1.07205507699324
$arr[2,15] = 1.007530323615285
$arr[2,16] = 1.014249021641849
$arr[2,17] = 1.007530323615285
$arr[2,18] = 0.9902605789396342
$arr[2,19] = 0.9758675133183052
$arr[2,20] = 0.9933665728410391
$arr[3,0] = 0.8477808095912102
$arr[3,1] = 1.267623904252883
$arr[3,2] = 0.9505008452391173
$arr[3,3] = 1.01298600756196
$arr[3,4] = 0.8477808095912102
$arr[3,5] = 0.9498403545749273
$arr[3,6] = 0.8732505114308362
$arr[3,7] = 0.9505008452391173
$arr[3,8] = 0.9505008452391173
$arr[3,9] = 1.128868980282419
$arr[3,10] = 0.9130528320821308
$arr[3,11] = 0.9505008452391173
$arr[3,12] = 1.267623904252883
$arr[3,13] = 1.057702356922047
$arr[3,14] = 1.090338368167507
$arr[3,15] = 1.021968519694403
$arr[3,16] = 1.009485848642075
$arr[3,17] = 1.021968519694403
$arr[3,18] = 0.9947395977913354
$arr[3,19] = 0.9858918472808916
$arr[3,20] = 0.9929880306269355
$arr[4,0] = 0.4385092129803458
$arr[4,1] = 2.334482157112217
$arr[4,2] = 0.5747055985020091
$arr[4,3] = 1.116024297814418
$arr[4,4] = 0.4385092129803458
$arr[4,5] = 0.6895102014881167
$arr[4,6] = 0.5166063997628573
$arr[4,7] = 0.5747055985020091
$arr[4,8] = 0.5747055985020091
$arr[4,9] = 1.558378770403783
$arr[4,10] = 0.6240961799865764
$arr[4,11] = 0.5747055985020091
$arr[4,12] = 2.334482157112217
$arr[4,13] = 1.386495685046282
$arr[4,14] = 1.479289168549397
$arr[4,15] = 1.115898989531524
$arr[4,16] = 1.132362516693046
$arr[4,17] = 1.115898989531524
$arr[4,18] = 0.9929482871452873
$arr[4,19] = 0.9092997494166315
$arr[4,20] = 0.9815391022562905
$arr[5,0] = 1.010403626782295
$arr[5,1] = 0.989305322038683
$arr[5,2] = 1.040420257879107
$arr[5,3] = 0.983625000792049
$arr[5,4] = 1.010403626782295
$arr[5,5] = 1.006274641959758
$arr[5,6] = 0.9966093621194336
$arr[5,7] = 1.040420257879107
$arr[5,8] = 1.040420257879107
$arr[5,9] = 1.002101409424883
$arr[5,10] = 0.9930398909249045
$arr[5,11] = 1.040420257879107
$arr[5,12] = 0.989305322038683
$arr[5,13] = 0.9998544744104889
$arr[5,14] = 0.9911726064817938
$arr[5,15] = 1.013376402233362
$arr[5,16] = 0.9975829465819608
$arr[5,17] = 1.013376402233362
$arr[5,18] = 1.008292274406247
$arr[5,19] = 1.01471787110082
$arr[5,20] = 1.002722438990139
$arr[6,0] = 1.005127200166364
$arr[6,1] = 0.9974414663356086
$arr[6,2] = 1.008239703266212
$arr[6,3] = 0.9980009742630671
$arr[6,4] = 1.005127200166364
$arr[6,5] = 0.9994180783567008
$arr[6,6] = 1.00023610656259
$arr[6,7] = 1.008239703266212
$arr[6,8] = 1.008239703266212
$arr[6,9] = 0.9999238758764856
$arr[6,10] = 0.9973096626796888
$arr[6,11] = 1.008239703266212
$arr[6,12] = 0.9974414663356086
$arr[6,13] = 1.001284333250986
$arr[6,14] = 0.9973755645076487
$arr[6,15] = 1.003602789922728
$arr[6,16] = 0.999959443060554
$arr[6,17] = 1.003602789922728
$arr[6,18] = 1.002029508111968
$arr[6,19] = 1.003271547142817
$arr[6,20] = 1.00071213343834
$arr[7,0] = 1.018061987890342
$arr[7,1] = 0.9845095129184515
$arr[7,2] = 1.066965948100107
$arr[7,3] = 0.9724853073651298
$arr[7,4] = 1.018061987890342
$arr[7,5] = 1.010594042399899
$arr[7,6] = 0.9942119729208772
$arr[7,7] = 1.066965948100107
$arr[7,8] = 1.066965948100107
$arr[7,9] = 1.003216506674783
$arr[7,10] = 0.9886094469464985
$arr[7,11] = 1.066965948100107
$arr[7,12] = 0.9845095129184515
$arr[7,13] = 1.001285750404397
$arr[7,14] = 0.9865594799324751
$arr[7,15] = 1.0231791496363
$arr[7,16] = 0.9970603159184307
$arr[7,17] = 1.0231791496363
$arr[7,18] = 1.01453672396385
$arr[7,19] = 1.025022568791101
$arr[7,20] = 1.004831840652011
$arr[8,0] = 0.4433765392047035
$arr[8,1] = 2.381196502040565
$arr[8,2] = 0.5598695232096317
$arr[8,3] = 1.118527375267104
$arr[8,4] = 0.4433765392047035
$arr[8,5] = 0.6719721326797238
$arr[8,6] = 0.5059101400744392
$arr[8,7] = 0.5598695232096317
$arr[8,8] = 0.5598695232096317
$arr[8,9] = 1.578326877268044
$arr[8,10] = 0.6048502953611316
$arr[8,11] = 0.5598695232096317
$arr[8,12] = 2.381196502040565
$arr[8,13] = 1.412286520622634
$arr[8,14] = 1.493023398700848
$arr[8,15] = 1.128147521484967
$arr[8,16] = 1.143141112202133
$arr[8,17] = 1.128147521484967
$arr[8,18] = 0.997323214954008
$arr[8,19] = 0.9098324766051327
$arr[8,20] = 0.9830036731381677
$arr[9,0] = 0.4684520838985147
$arr[9,1] = 1.29652235431551
$arr[9,2] = 1.193828020659943
$arr[9,3] = 0.9336955637600742
$arr[9,4] = 0.4684520838985147
$arr[9,5] = 1.07022058055835
$arr[9,6] = 0.6411936653680785
$arr[9,7] = 1.193828020659943
$arr[9,8] = 1.193828020659943
$arr[9,9] = 1.248727256418112
$arr[9,10] = 0.8600429324605693
$arr[9,11] = 1.193828020659943
$arr[9,12] = 1.29652235431551
$arr[9,13] = 0.8824872191070123
$arr[9,14] = 1.07828264338804
$arr[9,15] = 0.9862674862913227
$arr[9,16] = 0.8750057902248646
$arr[9,17] = 0.9862674862913227
$arr[9,18] = 0.9547113478336343
$arr[9,19] = 1.002534682398896
$arr[9,20] = 0.964085307179894
$arr[10,0] = 0.8192004398949468
$arr[10,1] = 1.404667167479996
$arr[10,2] = 0.8781127026111538
$arr[10,3] = 1.035454017002105
$arr[10,4] = 0.8192004398949468
$arr[10,5] = 0.9046973195894729
$arr[10,6] = 0.8388834451557873
$arr[10,7] = 0.8781127026111538
$arr[10,8] = 0.8781127026111538
$arr[10,9] = 1.181168457233681
$arr[10,10] = 0.8741875467368424
$arr[10,11] = 0.8781127026111538
$arr[10,12] = 1.404667167479996
$arr[10,13] = 1.111933803687472
$arr[10,14] = 1.139427357108419
$arr[10,15] = 1.033993436662032
$arr[10,16] = 1.032685051370595
$arr[10,17] = 1.033993436662032
$arr[10,18] = 0.9940419641807349
$arr[10,19] = 0.9708561118668186
$arr[10,20] = 0.9920463869629982
$arr[11,0] = 1.164790025901699
$arr[11,1] = 0.9865520076879692
$arr[11,2] = 0.9819036849905249
$arr[11,3] = 0.9964452724376196
$arr[11,4] = 1.164790025901699
$arr[11,5] = 0.9607708408790863
$arr[11,6] = 1.065374401430652
$arr[11,7] = 0.9819036849905249
$arr[11,8] = 0.9819036849905249
$arr[11,9] = 0.9782546329189276
$arr[11,10] = 0.9863329662155578
$arr[11,11] = 0.9819036849905249
$arr[11,12] = 0.9865520076879692
$arr[11,13] = 1.075671016794834
$arr[11,14] = 0.9864424869517635
$arr[11,15] = 1.044415239526731
$arr[11,16] = 1.045891666601742
$arr[11,17] = 1.044415239526731
$arr[11,18] = 1.029894671198938
$arr[11,19] = 1.020296473957255
$arr[11,20] = 1.015052979057755
$arr[12,0] = 0.02330075600000002
$arr[12,1] = 3.442862200000002
$arr[12,2] = 0.2138068799999999
$arr[12,3] = 1.209990899999999
$arr[12,4] = 0.02330075600000002
$arr[12,5] = 0.4216041200000005
$arr[12,6] = 0.1424155799999999
$arr[12,7] = 0.2138068799999999
$arr[12,8] = 0.2138068799999999
$arr[12,9] = 2.008686600000002
$arr[12,10] = 0.31618718
$arr[12,11] = 0.2138068799999999
$arr[12,12] = 3.442862200000002
$arr[12,13] = 1.733081478000001
$arr[12,14] = 1.879524690000001
$arr[12,15] = 1.226656612000001
$arr[12,16] = 1.260783378666667
$arr[12,17] = 1.226656612000001
$arr[12,18] = 0.9990392540000006
$arr[12,19] = 0.8419927792000005
$arr[12,20] = 0.9723567770000004
$arr[13,0] = 6.8875632
$arr[13,1] = 0.023300756
$arr[13,2] = 0.0017993233
$arr[13,3] = 0.7731529700000001
$arr[13,4] = 6.8875632
$arr[13,5] = 0.14230601
$arr[13,6] = 3.193909
$arr[13,7] = 0.0017993233
$arr[13,8] = 0.0017993233
$arr[13,9] = 0.014853478
$arr[13,10] = 0.95836551
$arr[13,11] = 0.0017993233
$arr[13,12] = 0.023300756
$arr[13,13] = 3.455431978
$arr[13,14] = 0.490833133
$arr[13,15] = 2.3042210931
$arr[13,16] = 2.623076488666667
$arr[13,17] = 2.3042210931
$arr[13,18] = 1.967757197325
$arr[13,19] = 1.57456562252
$arr[13,20] = 1.4994062809125
$arr[14,0] = 0.023300756
$arr[14,1] = 3.4428622
$arr[14,2] = 0.21380688
$arr[14,3] = 1.2099909
$arr[14,4] = 0.023300756
$arr[14,5] = 0.42160412
$arr[14,6] = 0.14241558
$arr[14,7] = 0.21380688
$arr[14,8] = 0.21380688
$arr[14,9] = 2.0086866
$arr[14,10] = 0.31618718
$arr[14,11] = 0.21380688
$arr[14,12] = 3.4428622
$arr[14,13] = 1.733081478
$arr[14,14] = 1.87952469
$arr[14,15] = 1.226656612
$arr[14,16] = 1.260783378666667
$arr[14,17] = 1.226656612
$arr[14,18] = 0.999039254
$arr[14,19] = 0.8419927792000002
$arr[14,20] = 0.9723567770000001
$arr[15,0] = 3.4457045
$arr[15,1] = 0.21658931
$arr[15,2] = 0.017301633
$arr[15,3] = 1.3509898
$arr[15,4] = 3.4457045
$arr[15,5] = 0.38954228
$arr[15,6] = 2.5451957
$arr[15,7] = 0.017301633
$arr[15,8] = 0.017301633
$arr[15,9] = 0.13660416
$arr[15,10] = 1.3306803
$arr[15,11] = 0.017301633
$arr[15,12] = 0.21658931
$arr[15,13] = 1.831146905
$arr[15,14] = 0.773634805
$arr[15,15] = 1.226531814333333
$arr[15,16] = 1.664324703333333
$arr[15,17] = 1.226531814333333
$arr[15,18] = 1.25256893575
$arr[15,19] = 1.0055154752
$arr[15,20] = 1.179075960375
$arr[16,0] = 1.302335194706849
$arr[16,1] = 0.7548764528219182
$arr[16,2] = 1.805248583839726
$arr[16,3] = 0.6620594544383563
$arr[16,4] = 1.302335194706849
$arr[16,5] = 1.127535245342466
$arr[16,6] = 0.9698589353424658
$arr[16,7] = 1.805248583839726
$arr[16,8] = 1.805248583839726
$arr[16,9] = 1.013239659095891
$arr[16,10] = 0.8724815289041099
$arr[16,11] = 1.805248583839726
$arr[16,12] = 0.7548764528219182
$arr[16,13] = 1.028605823764384
$arr[16,14] = 0.813678990863014
$arr[16,15] = 1.287486743789498
$arr[16,16] = 0.9765643921442925
$arr[16,17] = 1.287486743789498
$arr[16,18] = 1.183735440068151
$arr[16,19] = 1.308038068822466
$arr[16,20] = 1.063454381811473
$arr[17,0] = 2.313446413578948
$arr[17,1] = 1.341504682947369
$arr[17,2] = 0.08656499415263157
$arr[17,3] = 1.35322254368421
$arr[17,4] = 2.313446413578948
$arr[17,5] = 0.4148450784210526
$arr[17,6] = 1.598232091578947
$arr[17,7] = 0.08656499415263157
$arr[17,8] = 0.08656499415263157
$arr[17,9] = 0.7895791353684212
$arr[17,10] = 0.921983937368421
$arr[17,11] = 0.08656499415263157
$arr[17,12] = 1.341504682947369
$arr[17,13] = 1.827475548263158
$arr[17,14] = 1.131744310157895
$arr[17,15] = 1.247172030226316
$arr[17,16] = 1.525645011298246
$arr[17,17] = 1.247172030226316
$arr[17,18] = 1.165875007011842
$arr[17,19] = 0.95001300444
$arr[17,20] = 1.1024223596375
$arr[18,0] = 1.342166235863158
$arr[18,1] = 0.8159869238421055
$arr[18,2] = 1.739638943468421
$arr[18,3] = 0.6779404597894736
$arr[18,4] = 1.342166235863158
$arr[18,5] = 1.095994567894737
$arr[18,6] = 0.9773479268421053
$arr[18,7] = 1.739638943468421
$arr[18,8] = 1.739638943468421
$arr[18,9] = 1.026298970947369
$arr[18,10] = 0.8589722563157896
$arr[18,11] = 1.739638943468421
$arr[18,12] = 0.8159869238421055
$arr[18,13] = 1.079076579852632
$arr[18,14] = 0.8374795900789476
$arr[18,15] = 1.299264034391228
$arr[18,16] = 1.005708472007018
$arr[18,17] = 1.299264034391228
$arr[18,18] = 1.189191089872369
$arr[18,19] = 1.299280660591579
$arr[18,20] = 1.066793285620395
$arr[19,0] = 3.635942637510909
$arr[19,1] = 0.2404702404141349
$arr[19,2] = 0.04375825964634195
$arr[19,3] = 1.212490319680187
$arr[19,4] = 3.635942637510909
$arr[19,5] = 0.4758779848770468
$arr[19,6] = 2.492409552540417
$arr[19,7] = 0.04375825964634195
$arr[19,8] = 0.04375825964634195
$arr[19,9] = 0.177864771687463
$arr[19,10] = 1.346905311395828
$arr[19,11] = 0.04375825964634195
$arr[19,12] = 0.2404702404141349
$arr[19,13] = 1.938206438962522
$arr[19,14] = 0.7936877759049816
$arr[19,15] = 1.306723712523795
$arr[19,16] = 1.741106063106957
$arr[19,17] = 1.306723712523795
$arr[19,18] = 1.316769112241803
$arr[19,19] = 1.062166941722711
$arr[19,20] = 1.203214884719041
$arr[20,0] = 1.139827550582419
$arr[20,1] = 0.8386491562959589
$arr[20,2] = 0.3732002823117419
$arr[20,3] = 1.266482888926951
$arr[20,4] = 1.139827550582419
$arr[20,5] = 0.9237306816719273
$arr[20,6] = 1.333188373975093
$arr[20,7] = 0.3732002823117419
$arr[20,8] = 0.3732002823117419
$arr[20,9] = 0.7202106691882119
$arr[20,10] = 1.286060510114968
$arr[20,11] = 0.3732002823117419
$arr[20,12] = 0.8386491562959589
$arr[20,13] = 0.9892383534391891
$arr[20,14] = 1.062354833205463
$arr[20,15] = 0.78389232973004
$arr[20,16] = 1.088179072331115
$arr[20,17] = 0.7838923297300401
$arr[20,18] = 0.9094343748262721
$arr[20,19] = 0.802187556323366
$arr[20,20] = 0.9851687641334088
$arr[21,0] = 0.7496650844597287
$arr[21,1] = 0.7997830995241429
$arr[21,2] = 0.9597020113045045
$arr[21,3] = 1.019123428434768
$arr[21,4] = 0.7496650844597287
$arr[21,5] = 1.137698033442856
$arr[21,6] = 1.001317200607892
$arr[21,7] = 0.9597020113045045
$arr[21,8] = 0.9597020113045045
$arr[21,9] = 0.9149525429399392
$arr[21,10] = 1.153831219304493
$arr[21,11] = 0.9597020113045045
$arr[21,12] = 0.7997830995241429
$arr[21,13] = 0.7747240919919358
$arr[21,14] = 0.9768071594143179
$arr[21,15] = 0.8363833984294587
$arr[21,16] = 0.9010931344294549
$arr[21,17] = 0.8363833984294587
$arr[21,18] = 0.9157453536482172
$arr[21,19] = 0.9245366851794745
$arr[21,20] = 0.9670090775022904
$arr[22,0] = 1.268981835564235
$arr[22,1] = 1.212438411798979
$arr[22,2] = 1.142429010762399
$arr[22,3] = 0.978023070069021
$arr[22,4] = 1.268981835564235
$arr[22,5] = 0.830276611624532
$arr[22,6] = 0.986909305837664
$arr[22,7] = 1.142429010762399
$arr[22,8] = 1.142429010762399
$arr[22,9] = 1.098423037694314
$arr[22,10] = 0.7945051915643689
$arr[22,11] = 1.142429010762399
$arr[22,12] = 1.212438411798979
$arr[22,13] = 1.240710123681607
$arr[22,14] = 1.003471801681674
$arr[22,15] = 1.207949752708538
$arr[22,16] = 1.091975146309194
$arr[22,17] = 1.207949752708538
$arr[22,18] = 1.104588612422496
$arr[22,19] = 1.112156692090476
$arr[22,20] = 1.038998309364439
$arr[23,0] = 0.01541286061306674
$arr[23,1] = 0.5393055635557686
$arr[23,2] = 3.239353708501929
$arr[23,3] = 0.274247890737627
$arr[23,4] = 0.01541286061306674
$arr[23,5] = 1.636972449125893
$arr[23,6] = 0.1483227251941897
$arr[23,7] = 3.239353708501929
$arr[23,8] = 3.239353708501929
$arr[23,9] = 1.358578469365289
$arr[23,10] = 0.6052023913083295
$arr[23,11] = 3.239353708501929
$arr[23,12] = 0.5393055635557686
$arr[23,13] = 0.2773592120844177
$arr[23,14] = 0.5722539774320491
$arr[23,15] = 1.264690710890255
$arr[23,16] = 0.3866402718257216
$arr[23,17] = 1.264690710890255
$arr[23,18] = 1.099818630994773
$arr[23,19] = 1.527725646496204
$arr[23,20] = 0.9771745073002615
$arr[24,0] = 0.2282703272202107
$arr[24,1] = 1.409578004869539
$arr[24,2] = 0.6710635157005971
$arr[24,3] = 1.013506945489436
$arr[24,4] = 0.2282703272202107
$arr[24,5] = 1.131783605331437
$arr[24,6] = 0.6058989482421134
$arr[24,7] = 0.6710635157005971
$arr[24,8] = 0.6710635157005971
$arr[24,9] = 1.284452279092895
$arr[24,10] = 1.020287393186638
$arr[24,11] = 0.6710635157005971
$arr[24,12] = 1.409578004869539
$arr[24,13] = 0.8189241660448749
$arr[24,14] = 1.214932699028089
$arr[24,15] = 0.7696372825967823
$arr[24,16] = 0.886045241758796
$arr[24,17] = 0.7696372825967823
$arr[24,18] = 0.8322998102442463
$arr[24,19] = 0.8000525513355164
$arr[24,20] = 0.9206051273916083
$arr[25,0] = 0.8799281347641448
$arr[25,1] = 0.9664443216067362
$arr[25,2] = 1.063381027047922
$arr[25,3] = 0.9664400315857395
$arr[25,4] = 0.8799281347641448
$arr[25,5] = 1.039176766960464
$arr[25,6] = 0.998687515102231
$arr[25,7] = 1.063381027047922
$arr[25,8] = 1.063381027047922
$arr[25,9] = 0.9996533598188774
$arr[25,10] = 1.048176835544121
$arr[25,11] = 1.063381027047922
$arr[25,12] = 0.9664443216067362
$arr[25,13] = 0.9231862281854405
$arr[25,14] = 1.007310578575429
$arr[25,15] = 0.9699178278062677
$arr[25,16] = 0.9648497639716674
$arr[25,17] = 0.9699178278062677
$arr[25,18] = 0.9894825797407312
$arr[25,19] = 1.004262269202169
$arr[25,20] = 0.9952359990537796
$arr[26,0] = 1.127388372737343
$arr[26,1] = 0.850064924454385
$arr[26,2] = 1.13809338825602
$arr[26,3] = 0.8886509388033247
$arr[26,4] = 1.127388372737343
$arr[26,5] = 1.068205745052721
$arr[26,6] = 1.098889913839875
$arr[26,7] = 1.13809338825602
$arr[26,8] = 1.13809338825602
$arr[26,9] = 0.9308148223452857
$arr[26,10] = 1.084414291141361
$arr[26,11] = 1.13809338825602
$arr[26,12] = 0.850064924454385
$arr[26,13] = 0.9887266485958638
$arr[26,14] = 0.9672396077978732
$arr[26,15] = 1.038515561815916
$arr[26,16] = 1.020622529444363
$arr[26,17] = 1.038515561815916
$arr[26,18] = 1.049990244147277
$arr[26,19] = 1.067610872969026
$arr[26,20] = 1.02331529957879
$arr[27,0] = 1.671020857415538
$arr[27,1] = 0.8804038564634685
$arr[27,2] = 1.377010916546527
$arr[27,3] = 0.8335011395042078
$arr[27,4] = 1.671020857415538
$arr[27,5] = 0.8943972772203008
$arr[27,6] = 1.23937203971718
$arr[27,7] = 1.377010916546527
$arr[27,8] = 1.377010916546527
$arr[27,9] = 0.9100717848401919
$arr[27,10] = 0.9239339613005311
$arr[27,11] = 1.377010916546527
$arr[27,12] = 0.8804038564634685
$arr[27,13] = 1.275712356939503
$arr[27,14] = 0.9021689088819997
$arr[27,15] = 1.309478543475178
$arr[27,16] = 1.158452891726512
$arr[27,17] = 1.309478543475177
$arr[27,18] = 1.213092397931516
$arr[27,19] = 1.245876101654518
$arr[27,20] = 1.091213979125993
$ws.Range("C4:W31").Value = $arr

Write-Output "Edit complete"
